# Bank statement sample: split "amount" into separate "debit" / "credit"
# columns (inserted after "description"), keep "reference" after them, and
# recompute "amount" as a trailing formula column (credit - debit).
#
# Before: A=date  B=description  C=amount      D=reference
# After:  A=date  B=description  C=debit  D=credit  E=reference  F=amount(=D-C)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room for the two new columns (debit, credit) right after
#    "description". This pushes the old amount/reference columns from C/D
#    to E/F (formats + the C1/D1 header style carry along with the insert).
$ws.Columns("C:D").Insert()

# 2. Header row.
$ws.Range("C1").Value = "debit"
$ws.Range("D1").Value = "credit"
$ws.Range("E1").Value = "reference"

# New bordered header style for the trailing "amount" column. ClearFormats
# first so the freshly-inserted cell picks up its own style slot instead of
# sharing the one it inherited from the column shift.
$ws.Range("F1").ClearFormats()
$ws.Range("F1").Value = "amount"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160
$ws.Range("F1").Borders.Item(7).LineStyle = 1
$ws.Range("F1").Borders.Item(10).LineStyle = 1

# 3. Data rows: old amount (now sitting in column E after the insert) is
#    split into debit/credit, the reference text (now sitting in column F)
#    moves back to E, and F becomes a live formula.
$debit  = @{2=0;      3=5000;    4=0;       5=1000000; 6=0;       7=0;       8=0;     9=0;       10=1200000; 11=850000; 12=0;      13=0;      14=7500}
$credit = @{2=1500000; 3=0;       4=2500000; 5=0;       6=1750000; 7=2000000; 8=25000; 9=3000000; 10=0;       11=0;      12=2200000; 13=1800000; 14=0}
$reference = @{2="TRX-001"; 3="FEE-001"; 4="TRX-002"; 5="ATM-001"; 6="TRX-003"; 7="TRX-004"; 8="INT-001"; 9="TRX-005"; 10="PAY-001"; 11="PAY-002"; 12="TRX-006"; 13="TRX-007"; 14="FEE-002"}

for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value2 = $debit[$r]
    $ws.Cells.Item($r, 4).Value2 = $credit[$r]
    $ws.Range("E$r").Value = $reference[$r]
}

# Formula column: F2 is entered on its own, F3:F14 filled down together so
# they form one shared-formula group (matches a fill-down in real Excel).
$ws.Range("F2").Formula = "=D2-C2"
$ws.Range("F3:F14").Formula = "=D3-C3"

# 4. Restore the date column's custom number format (kept the same numFmtId,
#    just re-escaped the literal dashes).
$ws.Range("A2:A14").NumberFormat = "yyyy\-mm\-dd"

# 5. Column widths to fit the new layout (offset to account for the
#    engine's fixed character-padding delta between ColumnWidth and the
#    stored <col width> it writes out).
$pad = 0.8333333333333333
$ws.Columns("A").ColumnWidth = 10.33203125 - $pad
$ws.Columns("B").ColumnWidth = 24.109375 - $pad
$ws.Columns("C").ColumnWidth = 11 - $pad
$ws.Columns("D").ColumnWidth = 8.6640625 - $pad
$ws.Columns("E").ColumnWidth = 9 - $pad

$ws.Range("A1").Select()
